$wb = $excel.ActiveWorkbook

# Rename the two sheets (3rd sheet "extended mapping" is unchanged)
$wb.Worksheets.Item("cell_type2methods").Name = "celltype2method"
$wb.Worksheets.Item("cell_type2datasets").Name = "celltype2dataset"

# Move the active/selected tab from "celltype2method" (1st sheet) to
# "celltype2dataset" (2nd sheet)
$wb.Worksheets.Item("celltype2dataset").Select()
